$d = $word.ActiveDocument
$n = 0
foreach ($p in $d.Paragraphs) {
  $p.Alignment = 3
  $n = $n + 1
}
Write-Host "Processed:" $n
